$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.030.11'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '1.832.25'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("D4").Value = '''0.9990'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''242.06'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = '''0.6251'
$ws.Range("E6").Value = '  -5.69%  '
$ws.Range("D8").Value = '''0.07578'
$ws.Range("E8").Value = '  +1.61%  '
$ws.Range("D9").Value = '''0.2919'
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").Value = '''22.52'
$ws.Range("E10").Value = '  -3.20%  '
$ws.Range("D11").Value = '''0.07721'
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("D12").Value = '1.832.16'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = '''4.951'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").Value = '''0.6637'
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '''0.00001015'
$ws.Range("E15").Value = '  +16.34%  '
$ws.Range("D16").Value = '''82.60'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = '''6.022'
$ws.Range("E17").Value = '  -2.50%  '
$ws.Range("D18").Value = '29.005.02'
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("D19").Value = '''226.40'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("E20").Value = '  -1.54%  '
$ws.Range("D21").Value = '''0.9994'
$ws.Range("D22").Value = '''7.162'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '''158.03'
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").Value = '''8.466'
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("D26").Value = '''0.1372'
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = '''1.488'
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("D29").Value = '''4.089'
$ws.Range("D30").Value = '''4.011'
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = '''0.05196'
$ws.Range("E32").Value = '  -3.45%  '
$ws.Range("D33").Value = '''1.843'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").Value = '''0.7355'
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("E35").Value = '  -1.84%  '
$ws.Range("D36").Value = '''2.698'
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("D37").Value = '1.240.03'
$ws.Range("E37").Value = '  -4.63%  '
$ws.Range("D38").Value = '''2.754'
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("D40").Value = '''6.323'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("D41").Value = '''0.8942'
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("D43").Value = '''101.35'
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("D44").Value = '1.978.11'
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("E45").Value = '  +1.76%  '
$ws.Range("D46").Value = '''64.00'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("D47").Value = '''0.5107'
$ws.Range("D48").Value = '''0.4024'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").Value = '''8.844'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").Value = '''0.05753'
$ws.Range("E51").Value = '  -6.83%  '
